$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.003.85'
$ws.Range('E2').Value = '  -5.68%  '
$ws.Range('D3').Value = '2.681.31'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.42'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.71'
$ws.Range('E6').Value = '  -5.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.595'
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.00'
$ws.Range('E10').Value = '  -5.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.02'
$ws.Range('E12').Value = '  -3.25%  '
$ws.Range('D13').Value = '3.104.31'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.107'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').Value = '2.709.25'
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.925'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.08'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '45.148.47'
$ws.Range('E18').Value = '  -5.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.86'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000101'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.64'
$ws.Range('E21').Value = '  -4.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.07'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '278.28'
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.01'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '31.08'
$ws.Range('E26').Value = '  +1.71%  '
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.54'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.71'
$ws.Range('E30').Value = '  -5.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.17'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.77'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.35'
$ws.Range('E33').Value = '  +3.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '154.08'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0836'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.120'
$ws.Range('E37').Value = '  -5.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.25'
$ws.Range('E38').Value = '  +11.98%  '
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.00'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.60'
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0324'
$ws.Range('E42').Value = '  -3.37%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.146.66'
$ws.Range('E43').Value = '  -2.31%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.93'
$ws.Range('E44').Value = '  -7.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.29'
$ws.Range('E46').Value = '  -4.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.45'
$ws.Range('E47').Value = '  -5.04%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.957.10'
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '111.21'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.60'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.198'
$ws.Range('E51').Value = '  -2.45%  '
